# Updated cryptos list on Mon Oct 23 06:46:44 UTC 2023 with GitHub Actions
#
# Refreshes the coin Price (column D) and Volume(1h) (column E) figures for
# each ranked row, and fixes the ranking order for two coin pairs whose
# Coin/Link/Price/Volume data had swapped places (rows 39/40 and 43/44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as plain text (e.g. "1.00", "30.72", "0.0₆0115").
# Pre-format the cells we are about to touch as Text so Excel doesn't
# reinterpret numeric-looking values (losing trailing zeros, etc.) when we
# assign them below.
$priceCells = @(2,3,5,8,9,12,14,15,17,18,21,22,25,26,33,35,37,39,40,42,43,44,45,47,48,49,50,51)
foreach ($r in $priceCells) {
    $ws.Range("D$r").NumberFormat = '@'
}

# Row => updated Coin / Link / Price / Volume(1h) values.
# Only the fields actually present for a row are written; rows not listed
# here (e.g. row 29) are left untouched.
$updates = @(
    @{ Row = 2;  D = '30.745.59';  E = '  +2.07%  ' },
    @{ Row = 3;  D = '1.689.09';   E = '  +2.70%  ' },
    @{ Row = 4;  E = '  +0.11%  ' },
    @{ Row = 5;  D = '221.37';     E = '  +2.51%  ' },
    @{ Row = 6;  E = '  -0.08%  ' },
    @{ Row = 7;  E = '  +0.09%  ' },
    @{ Row = 8;  D = '30.72';      E = '  +4.55%  ' },
    @{ Row = 9;  D = '0.266';      E = '  +1.42%  ' },
    @{ Row = 10; E = '  +1.65%  ' },
    @{ Row = 11; E = '  -0.98%  ' },
    @{ Row = 12; D = '1.933.59';   E = '  +2.95%  ' },
    @{ Row = 13; E = '  +12.03%  ' },
    @{ Row = 14; D = '0.625';      E = '  +8.56%  ' },
    @{ Row = 15; D = '1.699.35';   E = '  +3.46%  ' },
    @{ Row = 16; E = '  +1.62%  ' },
    @{ Row = 17; D = '30.739.45';  E = '  +2.01%  ' },
    @{ Row = 18; D = '66.52';      E = '  +2.18%  ' },
    @{ Row = 19; E = '  -0.49%  ' },
    @{ Row = 20; E = '  +0.59%  ' },
    @{ Row = 21; D = '1.00';       E = '  +0.06%  ' },
    @{ Row = 22; D = '10.27';      E = '  +3.16%  ' },
    @{ Row = 23; E = '  +1.80%  ' },
    @{ Row = 24; E = '  +0.66%  ' },
    @{ Row = 25; D = '157.33';     E = '  -1.14%  ' },
    @{ Row = 26; D = '15.89';      E = '  +0.53%  ' },
    @{ Row = 27; E = '  -0.30%  ' },
    @{ Row = 28; E = '  +0.49%  ' },
    @{ Row = 30; E = '  +1.32%  ' },
    @{ Row = 31; E = '  +0.57%  ' },
    @{ Row = 32; E = '  +1.53%  ' },
    @{ Row = 33; D = '1.516.26';   E = '  +5.22%  ' },
    @{ Row = 34; E = '  +2.35%  ' },
    @{ Row = 35; D = '1.74';       E = '  +4.27%  ' },
    @{ Row = 36; E = '  -0.66%  ' },
    @{ Row = 37; D = '83.45';      E = '  +6.71%  ' },
    @{ Row = 38; E = '  +4.12%  ' },
    @{ Row = 39; B = 'MXToken';    C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx';     D = '2.75';  E = '  -4.45%  ' },
    @{ Row = 40; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx';     D = '0.585'; E = '  +4.19%  ' },
    @{ Row = 41; E = '  +1.57%  ' },
    @{ Row = 42; D = '0.850';      E = '  +0.17%  ' },
    @{ Row = 43; B = 'RenderToken'; C = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D = '2.01';   E = '  -1.66%  ' },
    @{ Row = 44; B = 'Kaspa';      C = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas';           D = '0.0505'; E = '  +1.28%  ' },
    @{ Row = 45; D = '1.03';       E = '  -0.68%  ' },
    @{ Row = 46; E = '  +0.06%  ' },
    @{ Row = 47; D = '51.84';      E = '  -7.18%  ' },
    @{ Row = 48; D = '1.825.92';   E = '  +2.26%  ' },
    @{ Row = 49; D = '5.44';       E = '  +0.83%  ' },
    @{ Row = 50; D = '94.94';      E = '  +4.81%  ' },
    @{ Row = 51; D = '0.0₆0115';   E = '  -0.08%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Range("B$r").Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Range("C$r").Value = $u.C }
    if ($u.ContainsKey('D')) { $ws.Range("D$r").Value = $u.D }
    if ($u.ContainsKey('E')) { $ws.Range("E$r").Value = $u.E }
}
